$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.451.03"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "2.200.57"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.45"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.58"
$ws.Range("E7").Value = "  -2.99%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -4.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.23"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0938"
$ws.Range("E11").Value = "  -4.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.84"
$ws.Range("E13").Value = "  -4.23%  "
$ws.Range("D14").Value = "2.529.72"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.11"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.822"
$ws.Range("E16").Value = "  -3.66%  "
$ws.Range("D17").Value = "2.195.08"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "41.425.82"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("E19").Value = "  -11.30%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.25"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.57"
$ws.Range("E22").Value = "  +5.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "227.64"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  -5.55%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.23"
$ws.Range("E26").Value = "  -6.96%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.30"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.32"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0789"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.72"
$ws.Range("E33").Value = "  +4.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.28"
$ws.Range("E34").Value = "  -7.26%  "
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("E36").Value = "  -9.84%  "
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.79"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.09"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.24"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.52"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.194"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.53"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.84"
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("E49").Value = "  -4.38%  "
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").Value = "2.406.54"
$ws.Range("E51").Value = "  -1.23%  "
